$d = $word.ActiveDocument

# --- Exercise 1, question 1: "RGBA" -> "ARGB" label, and value reorder ---
$q1 = $d.Content
$q1.Find.Execute("What is the color with the following values of RGBA:")
$q1.Text = "What is the color with the following values of ARGB:"

$v1 = $d.Range($q1.End, $d.Content.End)
$v1.Find.Execute("(255, 255, 0, 0)")
$v1.Text = "(0, 255, 255, 0)"

# --- Exercise 1, question 2: "RGBA" -> "ARGB" label, and value reorder ---
$q2 = $d.Range($v1.End, $d.Content.End)
$q2.Find.Execute("What is the color with the following values of RGBA:")
$q2.Text = "What is the color with the following values of ARGB:"

$v2 = $d.Range($q2.End, $d.Content.End)
$v2.Find.Execute("(0, 0, 255, 255)")
$v2.Text = "(255, 0, 255, 0)"

# --- Header: "Chapter 2 ... Images and Fonts ..." -> "Chapter 2 ... Images, Fonts and Colors ..." ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRng = $hdr.Range
$oldHeader = "Chapter 2                                                                                                 Images and Fonts                                                                                                                                           "
$newHeader = "Chapter 2                                                                                       Images, Fonts and Colors                                                                                                                                           "
$hdrRng.Find.Execute($oldHeader, $true, $false, $false, $false, $false, $true, 1, $false, $newHeader, 2)
